$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# Set the active/selected cell on this sheet to E5 (as seen when the sheet was last saved)
$ws.Activate()
$ws.Range("E5").Select()

# Add the new values in column E for rows 3 and 4
$ws.Range("E3").Value = 20
$ws.Range("E4").Value = 20
